$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 7959.5386
$ws.Range("I9").Value = 233.33333
$ws.Range("J9").Value = 14582
$ws.Range("K9").Value = 233.33333
$ws.Range("L9").Value = 14582
$ws.Range("M9").Value = -64.33332999999999
$ws.Range("N9").Value = -14920

# Row 16
$ws.Range("H16").Value = 6668.136
$ws.Range("I16").Value = 5566.6113
$ws.Range("J16").Value = 11625
$ws.Range("K16").Value = 5566.6113
$ws.Range("L16").Value = 11625
$ws.Range("M16").Value = -5336.6113
$ws.Range("N16").Value = -12085

# Row 28
$ws.Range("H28").Value = 1093.875
$ws.Range("I28").Value = 1066.8
$ws.Range("K28").Value = 1066.8
$ws.Range("M28").Value = -581.8

# Row 32
$ws.Range("H32").Value = 4738.4375
$ws.Range("I32").Value = 5194.9165
$ws.Range("J32").Value = 4464.55
$ws.Range("K32").Value = 5194.9165
$ws.Range("L32").Value = 4464.55
$ws.Range("M32").Value = -4868.9165
$ws.Range("N32").Value = -5116.55

# Row 33
$ws.Range("H33").Value = 16403510
$ws.Range("I33").Value = 22963862
$ws.Range("K33").Value = 22963862
$ws.Range("M33").Value = -22963633

# Row 62
$ws.Range("H62").Value = 1811.9
$ws.Range("I62").Value = 1811.9
$ws.Range("K62").Value = 1811.9
$ws.Range("M62").Value = -1187.9

# Row 65
$ws.Range("H65").Value = 1811.9
$ws.Range("I65").Value = 1811.9
$ws.Range("K65").Value = 9059.5
$ws.Range("M65").Value = -5939.5

# Row 69
$ws.Range("H69").Value = 13004.5
$ws.Range("J69").Value = 10015
$ws.Range("L69").Value = 30045
$ws.Range("N69").Value = -31793

# Row 72
$ws.Range("H72").Value = 13004.5
$ws.Range("J72").Value = 10015
$ws.Range("L72").Value = 90135
$ws.Range("N72").Value = -98871

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Row 88
$ws.Range("H88").Value = 1935.5416
$ws.Range("I88").Value = 1590.75
$ws.Range("J88").Value = 2004.5
$ws.Range("K88").Value = 1590.75
$ws.Range("L88").Value = 2004.5
$ws.Range("M88").Value = -1184.75
$ws.Range("N88").Value = -2816.5

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Row 91
$ws.Range("H91").Value = 1935.5416
$ws.Range("I91").Value = 1590.75
$ws.Range("J91").Value = 2004.5
$ws.Range("K91").Value = 1590.75
$ws.Range("L91").Value = 2004.5
$ws.Range("M91").Value = -186.75
$ws.Range("N91").Value = -4812.5

# Row 98
$ws.Range("H98").Value = 838.08
$ws.Range("I98").Value = 831.5
$ws.Range("J98").Value = 996
$ws.Range("K98").Value = 831.5
$ws.Range("L98").Value = 996
$ws.Range("M98").Value = 666.5
$ws.Range("N98").Value = -3992

# Row 100
$ws.Range("H100").Value = 3106.3635
$ws.Range("J100").Value = 3218.75
$ws.Range("L100").Value = 3218.75
$ws.Range("N100").Value = -4300.75

# Row 107
$ws.Range("H107").Value = 1052.3667
$ws.Range("I107").Value = 1153.8695
$ws.Range("J107").Value = 718.8570999999999
$ws.Range("K107").Value = 1153.8695
$ws.Range("L107").Value = 718.8570999999999
$ws.Range("M107").Value = 766.1305
$ws.Range("N107").Value = -4558.8571

# Row 111
$ws.Range("H111").Value = 2176.5
$ws.Range("I111").Value = 1997.8
$ws.Range("K111").Value = 5993.4
$ws.Range("M111").Value = -2926.4

# Row 112
$ws.Range("H112").Value = 1633.909
$ws.Range("J112").Value = 1741.75
$ws.Range("L112").Value = 5225.25
$ws.Range("N112").Value = -7441.25

# Row 122
$ws.Range("H122").Value = 838.08
$ws.Range("I122").Value = 831.5
$ws.Range("J122").Value = 996
$ws.Range("K122").Value = 2494.5
$ws.Range("L122").Value = 2988
$ws.Range("M122").Value = -44.5
$ws.Range("N122").Value = -7888

# Row 131
$ws.Range("H131").Value = 7822.2144
$ws.Range("I131").Value = 3835
$ws.Range("J131").Value = 14999.2
$ws.Range("K131").Value = 11505
$ws.Range("L131").Value = 44997.60000000001
$ws.Range("M131").Value = -6465
$ws.Range("N131").Value = -55077.60000000001

# Row 135
$ws.Range("H135").Value = 980.23334
$ws.Range("I135").Value = 1022.1111
$ws.Range("K135").Value = 9198.999899999999
$ws.Range("M135").Value = -6663.999899999999

# Row 137
$ws.Range("H137").Value = 1425.7778
$ws.Range("I137").Value = 1272.8125
$ws.Range("J137").Value = 2649.5
$ws.Range("K137").Value = 3818.4375
$ws.Range("L137").Value = 7948.5
$ws.Range("M137").Value = -1268.4375
$ws.Range("N137").Value = -13048.5

# Row 138
$ws.Range("H138").Value = 4428.75
$ws.Range("I138").Value = 3949.5
$ws.Range("J138").Value = 4482
$ws.Range("K138").Value = 11848.5
$ws.Range("L138").Value = 13446
$ws.Range("M138").Value = -6708.5
$ws.Range("N138").Value = -23726

# Row 141
$ws.Range("H141").Value = 3230.6
$ws.Range("I141").Value = 3201.2778
$ws.Range("K141").Value = 9603.8334
$ws.Range("M141").Value = -4423.8334


# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 3331.7937
$ws.Range("I32").Value = 2436.1091
$ws.Range("K32").Value = 2436.1091
$ws.Range("M32").Value = -2149.1091

# Row 43
$ws.Range("H43").Value = 20766.834
$ws.Range("J43").Value = 20766.834
$ws.Range("L43").Value = 20766.834
$ws.Range("N43").Value = -21392.834

# Row 61
$ws.Range("H61").Value = 3319.2563
$ws.Range("I61").Value = 2782.12
$ws.Range("J61").Value = 4278.4287
$ws.Range("K61").Value = 2782.12
$ws.Range("L61").Value = 4278.4287
$ws.Range("M61").Value = -2570.12
$ws.Range("N61").Value = -4702.4287

# Row 74
$ws.Range("H74").Value = 3611.2104
$ws.Range("I74").Value = 2245.1
$ws.Range("J74").Value = 5129.1113
$ws.Range("K74").Value = 2245.1
$ws.Range("L74").Value = 5129.1113
$ws.Range("M74").Value = -1371.1
$ws.Range("N74").Value = -6877.1113

# Row 77
$ws.Range("H77").Value = 3611.2104
$ws.Range("I77").Value = 2245.1
$ws.Range("J77").Value = 5129.1113
$ws.Range("K77").Value = 11225.5
$ws.Range("L77").Value = 25645.5565
$ws.Range("M77").Value = -6857.5
$ws.Range("N77").Value = -34381.5565

# Row 102
$ws.Range("H102").Value = 911077.75
$ws.Range("I102").Value = 1001997.8
$ws.Range("J102").Value = 1877
$ws.Range("K102").Value = 1001997.8
$ws.Range("L102").Value = 1877
$ws.Range("M102").Value = -1000375.8
$ws.Range("N102").Value = -5121

# Row 122
$ws.Range("H122").Value = 19608896
$ws.Range("I122").Value = 930
$ws.Range("J122").Value = 83334780
$ws.Range("K122").Value = 2790
$ws.Range("L122").Value = 250004340
$ws.Range("M122").Value = -340
$ws.Range("N122").Value = -250009240

# Row 132
$ws.Range("H132").Value = 3470.5908
$ws.Range("I132").Value = 3404.647
$ws.Range("K132").Value = 10213.941
$ws.Range("M132").Value = -7683.940999999999

# Row 136
$ws.Range("H136").Value = 3319.2563
$ws.Range("I136").Value = 2782.12
$ws.Range("J136").Value = 4278.4287
$ws.Range("K136").Value = 8346.360000000001
$ws.Range("L136").Value = 12835.2861
$ws.Range("M136").Value = -5796.360000000001
$ws.Range("N136").Value = -17935.2861


# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 2446.7
$ws.Range("I20").Value = 2454
$ws.Range("K20").Value = 2454
$ws.Range("M20").Value = -2207

# Row 81
$ws.Range("H81").Value = 25243.857
$ws.Range("J81").Value = 25243.857
$ws.Range("L81").Value = 25243.857
$ws.Range("N81").Value = -27365.857

# Row 84
$ws.Range("H84").Value = 25243.857
$ws.Range("J84").Value = 25243.857
$ws.Range("L84").Value = 75731.571
$ws.Range("N84").Value = -86339.571

# Row 86
$ws.Range("H86").Value = 100002060
$ws.Range("I86").Value = 200001150
$ws.Range("J86").Value = 2981.6
$ws.Range("K86").Value = 200001150
$ws.Range("L86").Value = 2981.6
$ws.Range("M86").Value = -200000027
$ws.Range("N86").Value = -5227.6

# Row 89
$ws.Range("H89").Value = 100002060
$ws.Range("I89").Value = 200001150
$ws.Range("J89").Value = 2981.6
$ws.Range("K89").Value = 1000005750
$ws.Range("L89").Value = 14908
$ws.Range("M89").Value = -1000000134
$ws.Range("N89").Value = -26140

# Row 107
$ws.Range("H107").Value = 11747.6
$ws.Range("I107").Value = 14009.5
$ws.Range("K107").Value = 14009.5
$ws.Range("M107").Value = -12089.5


# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 6
$ws.Range("H6").Value = 2013199.6
$ws.Range("I6").Value = 2013199.6
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2013199.6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2013086.6
$ws.Range("N6").ClearContents()

# Row 7
$ws.Range("H7").Value = 206.57143
$ws.Range("I7").Value = 116.8125
$ws.Range("K7").Value = 116.8125
$ws.Range("M7").Value = -3.8125

# Row 31
$ws.Range("H31").Value = 1518.75
$ws.Range("I31").Value = 1220.7097
$ws.Range("J31").Value = 2545.3333
$ws.Range("K31").Value = 1220.7097
$ws.Range("L31").Value = 2545.3333
$ws.Range("M31").Value = -925.7097000000001
$ws.Range("N31").Value = -3135.3333

# Row 34
$ws.Range("H34").Value = 1518.75
$ws.Range("I34").Value = 1220.7097
$ws.Range("J34").Value = 2545.3333
$ws.Range("K34").Value = 1220.7097
$ws.Range("L34").Value = 2545.3333
$ws.Range("M34").Value = -1018.7097
$ws.Range("N34").Value = -2949.3333

# Row 51
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472
$ws.Range("M51").ClearContents()

# Row 58
$ws.Range("H58").Value = 100004170
$ws.Range("I58").Value = 66670740
$ws.Range("J58").Value = 125004240
$ws.Range("K58").Value = 66670740
$ws.Range("L58").Value = 125004240
$ws.Range("M58").Value = -66670537
$ws.Range("N58").Value = -125004646

# Row 60
$ws.Range("H60").Value = 43679.633
$ws.Range("J60").Value = 43679.633
$ws.Range("L60").Value = 43679.633
$ws.Range("N60").Value = -44701.633

# Row 61
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696
$ws.Range("M61").ClearContents()

# Row 105
$ws.Range("H105").Value = 1377.5454
$ws.Range("I105").Value = 1445.5
$ws.Range("K105").Value = 1445.5
$ws.Range("M105").Value = 301.5

# Row 132
$ws.Range("H132").Value = 8182.5
$ws.Range("I132").Value = 5780
$ws.Range("K132").Value = 17340
$ws.Range("M132").Value = -14810

# Row 134
$ws.Range("H134").Value = 1875.2258
$ws.Range("I134").Value = 1881.5358
$ws.Range("J134").Value = 1816.3334
$ws.Range("K134").Value = 5644.607400000001
$ws.Range("L134").Value = 5449.0002
$ws.Range("M134").Value = -3109.607400000001
$ws.Range("N134").Value = -10519.0002

# Row 136
$ws.Range("H136").Value = 100004170
$ws.Range("I136").Value = 66670740
$ws.Range("J136").Value = 125004240
$ws.Range("K136").Value = 200012220
$ws.Range("L136").Value = 375012720
$ws.Range("M136").Value = -200009670
$ws.Range("N136").Value = -375017820

# Row 141
$ws.Range("H141").Value = 28847.637
$ws.Range("I141").Value = 23499
$ws.Range("J141").Value = 30036.223
$ws.Range("K141").Value = 23499
$ws.Range("L141").Value = 30036.223
$ws.Range("M141").Value = -18319
$ws.Range("N141").Value = -40396.223


# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 7701.1665
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 4
$ws.Range("H4").Value = 387035.34
$ws.Range("I4").Value = 601859.9399999999
$ws.Range("K4").Value = 1805579.82
$ws.Range("M4").Value = -1805467.82

# Row 5
$ws.Range("H5").Value = 2061.4285
$ws.Range("I5").Value = 771.3333
$ws.Range("J5").Value = 3029
$ws.Range("K5").Value = 2313.9999
$ws.Range("L5").Value = 9087
$ws.Range("M5").Value = -2201.9999
$ws.Range("N5").Value = -9311

# Row 46
$ws.Range("H46").Value = 1350
$ws.Range("I46").Value = 648
$ws.Range("J46").Value = 1490.4
$ws.Range("K46").Value = 1944
$ws.Range("L46").Value = 4471.200000000001
$ws.Range("N46").Value = -4653.200000000001
$ws.Range("M46").Value = -1853

# Row 63
$ws.Range("H63").Value = 2212.6667
$ws.Range("J63").Value = 2002.3334
$ws.Range("L63").Value = 6007.0002
$ws.Range("N63").Value = -7505.0002

# Row 66
$ws.Range("H66").Value = 2212.6667
$ws.Range("J66").Value = 2002.3334
$ws.Range("L66").Value = 18021.0006
$ws.Range("N66").Value = -25509.0006

# Row 86
$ws.Range("H86").Value = 978.5294
$ws.Range("I86").Value = 402.7143
$ws.Range("J86").Value = 3665.6667
$ws.Range("K86").Value = 1208.1429
$ws.Range("L86").Value = 10997.0001
$ws.Range("M86").Value = -22.14289999999983
$ws.Range("N86").Value = -13369.0001

# Row 89
$ws.Range("H89").Value = 978.5294
$ws.Range("I89").Value = 402.7143
$ws.Range("J89").Value = 3665.6667
$ws.Range("K89").Value = 3624.4287
$ws.Range("L89").Value = 32991.0003
$ws.Range("M89").Value = 2303.5713
$ws.Range("N89").Value = -44847.0003

# Row 107
$ws.Range("H107").Value = 1508.6296
$ws.Range("I107").Value = 603.6667
$ws.Range("J107").Value = 1961.1111
$ws.Range("K107").Value = 1811.0001
$ws.Range("L107").Value = 5883.3333
$ws.Range("M107").Value = 108.9999
$ws.Range("N107").Value = -9723.3333

# Row 109
$ws.Range("H109").Value = 7388.9473
$ws.Range("I109").Value = 2774
$ws.Range("K109").Value = 8322
$ws.Range("M109").Value = -7282

# Row 127
$ws.Range("H127").Value = 57485.5
$ws.Range("J127").Value = 57485.5
$ws.Range("L127").Value = 172456.5
$ws.Range("N127").Value = -182376.5

# Row 130
$ws.Range("H130").Value = 6137.4
$ws.Range("I130").Value = 4993.5
$ws.Range("J130").Value = 6900
$ws.Range("K130").Value = 14980.5
$ws.Range("L130").Value = 20700
$ws.Range("M130").Value = -9960.5
$ws.Range("N130").Value = -30740

# Row 135
$ws.Range("H135").Value = 2061.4285
$ws.Range("I135").Value = 771.3333
$ws.Range("J135").Value = 3029
$ws.Range("K135").Value = 6941.9997
$ws.Range("L135").Value = 27261
$ws.Range("M135").Value = -4406.9997
$ws.Range("N135").Value = -32331

# Row 137
$ws.Range("H137").Value = 1989.3492
$ws.Range("J137").Value = 2076.5535
$ws.Range("L137").Value = 6229.6605
$ws.Range("N137").Value = -16429.6605

# Row 140
$ws.Range("H140").Value = 3054.1667
$ws.Range("I140").Value = 3009.6667
$ws.Range("K140").Value = 9029.000100000001
$ws.Range("M140").Value = -3849.000100000001


# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 2
$ws.Range("H2").Value = 1217.6666
$ws.Range("I2").Value = 52
$ws.Range("J2").Value = 2674.75
$ws.Range("K2").Value = 52
$ws.Range("L2").Value = 2674.75
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = -2900.75

# Row 63
$ws.Range("H63").Value = 31980
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 31980
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 31980
$ws.Range("N63").Value = -33352
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 31980
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 31980
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 95940
$ws.Range("N66").Value = -102804
$ws.Range("M66").ClearContents()

# Row 80
$ws.Range("H80").Value = 3465.4443
$ws.Range("I80").Value = 2674.5
$ws.Range("K80").Value = 2674.5
$ws.Range("M80").Value = -1676.5

# Row 83
$ws.Range("H83").Value = 3465.4443
$ws.Range("I83").Value = 2674.5
$ws.Range("K83").Value = 13372.5
$ws.Range("M83").Value = -8380.5

# Row 87
$ws.Range("H87").Value = 47966.668
$ws.Range("J87").Value = 47966.668
$ws.Range("L87").Value = 47966.668
$ws.Range("N87").Value = -50462.668

# Row 90
$ws.Range("H90").Value = 47966.668
$ws.Range("J90").Value = 47966.668
$ws.Range("L90").Value = 143900.004
$ws.Range("N90").Value = -156380.004

# Row 97
$ws.Range("H97").Value = 34484344
$ws.Range("I97").Value = 43479590
$ws.Range("K97").Value = 43479590
$ws.Range("M97").Value = -43479094

# Row 102
$ws.Range("H102").Value = 13337253
$ws.Range("I102").Value = 17337444
$ws.Range("J102").Value = 3279.5557
$ws.Range("K102").Value = 17337444
$ws.Range("L102").Value = 3279.5557
$ws.Range("M102").Value = -17335822
$ws.Range("N102").Value = -6523.5557

# Row 113
$ws.Range("H113").Value = 90913390
$ws.Range("I113").Value = 111115340
$ws.Range("J113").Value = 4606.5
$ws.Range("K113").Value = 111115340
$ws.Range("L113").Value = 4606.5
$ws.Range("M113").Value = -111113170
$ws.Range("N113").Value = -8946.5

# Row 122
$ws.Range("H122").Value = 62779.367
$ws.Range("I122").Value = 73061
$ws.Range("J122").Value = 24223.25
$ws.Range("K122").Value = 219183
$ws.Range("L122").Value = 72669.75
$ws.Range("M122").Value = -216733
$ws.Range("N122").Value = -77569.75

# Row 126
$ws.Range("H126").Value = 27781798
$ws.Range("I126").Value = 45457440
$ws.Range("K126").Value = 136372320
$ws.Range("M126").Value = -136369850

# Row 132
$ws.Range("H132").Value = 4928.6665
$ws.Range("I132").Value = 4915.3
$ws.Range("J132").Value = 4995.5
$ws.Range("K132").Value = 14745.9
$ws.Range("L132").Value = 14986.5
$ws.Range("M132").Value = -12215.9
$ws.Range("N132").Value = -20046.5

# Row 135
$ws.Range("H135").Value = 79394.71000000001
$ws.Range("J135").Value = 79394.71000000001
$ws.Range("L135").Value = 79394.71000000001
$ws.Range("N135").Value = -89534.71000000001


# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 3072.15
$ws.Range("I46").Value = 809.6
$ws.Range("J46").Value = 5334.7
$ws.Range("K46").Value = 809.6
$ws.Range("L46").Value = 5334.7
$ws.Range("M46").Value = -621.6
$ws.Range("N46").Value = -5710.7

# Row 68
$ws.Range("H68").Value = 6485.143
$ws.Range("J68").Value = 11000
$ws.Range("L68").Value = 11000
$ws.Range("N68").Value = -12498

# Row 71
$ws.Range("H71").Value = 6485.143
$ws.Range("J71").Value = 11000
$ws.Range("L71").Value = 55000
$ws.Range("N71").Value = -62488

# Row 125
$ws.Range("H125").Value = 69712.71000000001
$ws.Range("J125").Value = 69712.71000000001
$ws.Range("L125").Value = 69712.71000000001
$ws.Range("N125").Value = -79552.71000000001

# Row 132
$ws.Range("H132").Value = 12213.409
$ws.Range("I132").Value = 11039.8
$ws.Range("J132").Value = 23949.5
$ws.Range("K132").Value = 33119.39999999999
$ws.Range("L132").Value = 71848.5
$ws.Range("M132").Value = -30589.39999999999
$ws.Range("N132").Value = -76908.5

# Row 136
$ws.Range("H136").Value = 7410708.5
$ws.Range("I136").Value = 2701.8262
$ws.Range("K136").Value = 8105.4786
$ws.Range("M136").Value = -5555.4786

# Row 139
$ws.Range("H139").Value = 89703.75
$ws.Range("J139").Value = 89703.75
$ws.Range("L139").Value = 89703.75
$ws.Range("N139").Value = -99983.75


# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 3
$ws.Range("H3").Value = 4100
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 5150
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 5150
$ws.Range("M3").Value = -1886
$ws.Range("N3").Value = -5378

# Row 6
$ws.Range("H6").Value = 2339.6
$ws.Range("I6").Value = 1899
$ws.Range("J6").Value = 2449.75
$ws.Range("K6").Value = 1899
$ws.Range("L6").Value = 2449.75
$ws.Range("N6").Value = -2679.75
$ws.Range("M6").Value = -1784

# Row 62
$ws.Range("H62").Value = 9787.632
$ws.Range("I62").Value = 4412
$ws.Range("J62").Value = 12268.692
$ws.Range("K62").Value = 4412
$ws.Range("L62").Value = 12268.692
$ws.Range("M62").Value = -3788
$ws.Range("N62").Value = -13516.692

# Row 64
$ws.Range("H64").Value = 58666.668
$ws.Range("J64").Value = 58666.668
$ws.Range("L64").Value = 58666.668
$ws.Range("N64").Value = -59162.668

# Row 65
$ws.Range("H65").Value = 9787.632
$ws.Range("I65").Value = 4412
$ws.Range("J65").Value = 12268.692
$ws.Range("K65").Value = 22060
$ws.Range("L65").Value = 61343.45999999999
$ws.Range("M65").Value = -18940
$ws.Range("N65").Value = -67583.45999999999

# Row 67
$ws.Range("H67").Value = 58666.668
$ws.Range("J67").Value = 58666.668
$ws.Range("L67").Value = 58666.668
$ws.Range("N67").Value = -60382.668

# Row 113
$ws.Range("H113").Value = 1687.1
$ws.Range("I113").Value = 1430.1111
$ws.Range("K113").Value = 4290.3333
$ws.Range("M113").Value = -2120.3333

# Row 126
$ws.Range("H126").Value = 2385.6155
$ws.Range("I126").Value = 2234
$ws.Range("J126").Value = 2726.75
$ws.Range("K126").Value = 6702
$ws.Range("L126").Value = 8180.25
$ws.Range("M126").Value = -4232
$ws.Range("N126").Value = -13120.25

# Row 132
$ws.Range("H132").Value = 3154.6667
$ws.Range("I132").Value = 3162.15
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 9486.450000000001
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -6956.450000000001
$ws.Range("N132").Value = -14075

# Row 136
$ws.Range("H136").Value = 1811.6552
$ws.Range("I136").Value = 1621.56
$ws.Range("K136").Value = 4864.68
$ws.Range("M136").Value = -2314.68

# Row 139
$ws.Range("H139").Value = 69039.75999999999
$ws.Range("J139").Value = 69039.75999999999
$ws.Range("L139").Value = 69039.75999999999
$ws.Range("N139").Value = -79319.75999999999

